$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-13 Thursday" "2025-03-14 Friday"

Replace-Text "43×41=" "50×34="
Replace-Text "40×45=" "47×93="
Replace-Text "82×91=" "52×23="
Replace-Text "74×73=" "47×26="
Replace-Text "29×59=" "63×82="

Replace-Text "14×61=" "24×16="
Replace-Text "53×98=" "94×30="
Replace-Text "31×28=" "94×12="
Replace-Text "95×98=" "74×58="
Replace-Text "72×68=" "28×50="

Replace-Text "84×94=" "80×37="
Replace-Text "97×64=" "63×86="
Replace-Text "62×71=" "54×21="
Replace-Text "30×78=" "53×40="
Replace-Text "11×60=" "54×80="

Replace-Text "71×53=" "90×88="
Replace-Text "89×37=" "14×34="
Replace-Text "50×71=" "12×86="
Replace-Text "79×73=" "18×72="
Replace-Text "85×17=" "21×63="

Replace-Text "65×88=" "43×45="
Replace-Text "12×94=" "47×34="
Replace-Text "80×23=" "97×37="
Replace-Text "68×48=" "31×17="
Replace-Text "60×91=" "33×21="
